$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Gestão"
$ws.Range("C3").Value = "Desenho Técnico"
$ws.Range("C4").Value = "-"
$ws.Range("E4").Value = "-"
$ws.Range("F4").Value = "Circuitos Elétricos"
$ws.Range("B6").Value = "-"
